# Edit script: apply midterm question bank changes via Word COM interop
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: add a collapsed "_GoBack" bookmark right before the run that
# starts "During that lab, you will not be allowed ..." (2nd bullet point).
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("During that lab, you will not be allowed to use internet resources", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rng1)

# ---------------------------------------------------------------------------
# Change 2: question 3 prompt -- "following braces are balanced: (); []; {}."
# becomes "provided braces are balanced." (split across three runs).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("following braces are balanced", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $rng2.Paragraphs(1).Range
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="46DB3187" w14:textId="1167EB46" w:rsidR="00A53BA7" w:rsidRDefault="00A53BA7" w:rsidP="00A53BA7"><w:pPr><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t xml:space="preserve">3. </w:t></w:r><w:r w:rsidR="004B2683"><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t>Use a hash table to w</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t xml:space="preserve">rite a function called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t>isBalanced</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> that determines whether or not the </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t>provided braces are balanced</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t>.  E.g. "(ab[cd])" is balanced, "(ab" is not balanced, "(ab[cd)e]" is not balanced.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 3: isBalanced() signature -- "string text" becomes "string to_evaluate"
# and "unordered_map<string, string> parens" becomes "... braces"; the
# "_GoBack" bookmark that used to sit in this paragraph is gone (it moved to
# change 1 above).
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("bool isBalanced(string text", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para3 = $rng3.Paragraphs(1).Range
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="77C7FC8A" w14:textId="5D0821FD" w:rsidR="00A53BA7" w:rsidRDefault="00A53BA7" w:rsidP="00A53BA7"><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t>bool</w:t></w:r><w:r w:rsidRPr="00A53BA7"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t>isBalanced</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00A53BA7"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve">string </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t>to_evaluate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00B5132A"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00B5132A"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t>unordered_map</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00B5132A"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve">&lt;string, string&gt; </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t>braces</w:t></w:r><w:r w:rsidRPr="00A53BA7"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Courier New"/></w:rPr><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para3.InsertXML($xml3)
